# "Add files via upload" — progress-tracking workbook update.
# Folha1 (the only sheet) gets a new note in E5, row 18 (func #17) and
# row 21 (func #20) are marked "incorrect" (yellow fill) with updated
# values, and the long error note on row 18 is replaced with a short one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18 ("17;ConvertToBW_Otsu") now flagged incorrect: yellow fill on
# B:D, coreDiff/borderDiff replaced with the placeholder "~0.00", and the
# long diagnostic note replaced with a short one.
$ws.Range("B22:D22").Interior.Color = 65535
$ws.Range("C22").Value = "~0.00"
$ws.Range("D22").Value = "~0.00"
$ws.Range("E22").Value = "quase bom"

# Row 21 ("20;Roberts") now flagged incorrect too: yellow fill on B:D,
# coreDiff/borderDiff filled in with 0.
$ws.Range("B25:D25").Interior.Color = 65535
$ws.Range("C25").Value = 0
$ws.Range("D25").Value = 0

# New remark on row 5 ("1;Negative") — previously empty.
$ws.Range("E5").Value = "já saca a width do tabuleiro (nível 1)"

# Selection moved to the newly-edited cell E5 (was G15).
$ws.Range("E5").Select()
